$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuilt data for rows 309-433: a new SRB (Serbia) block is inserted
# (years 2006-2018) before TJK, and each of TJK/TKM/UKR/UZB's 1991-2018 series is
# split into a 1991-2005 chunk followed by a 2006-2018 chunk.
$data = @(
    @(2006, "SRB", 9, 1, 0, 8, 8),
    @(2007, "SRB", 9, 1, 1, 8, 8),
    @(2008, "SRB", 9, 1, 2, 8, 8),
    @(2009, "SRB", 9, 1, 3, 8, 8),
    @(2010, "SRB", 9, 1, 4, 8, 8),
    @(2011, "SRB", 9, 1, 5, 8, 8),
    @(2012, "SRB", 9, 1, 6, 8, 8),
    @(2013, "SRB", 9, 1, 7, 8, 8),
    @(2014, "SRB", 9, 1, 8, 8, 8),
    @(2015, "SRB", 9, 1, 9, 8, 8),
    @(2016, "SRB", 9, 1, 10, 8, 8),
    @(2017, "SRB", 9, 1, 11, 8, 8),
    @(2018, "SRB", 9, 1, 12, 8, 8),
    @(1991, "TJK", 1, 3, 0, -2, -2),
    @(1992, "TJK", 0, 6, 0, -6, -6),
    @(1993, "TJK", 0, 6, 1, -6, -6),
    @(1994, "TJK", 0, 6, 2, -6, -6),
    @(1995, "TJK", 0, 6, 3, -6, -6),
    @(1996, "TJK", 0, 6, 4, -6, -6),
    @(1997, "TJK", 0, 5, 0, -5, -5),
    @(1998, "TJK", 2, 3, 0, -1, -1),
    @(1999, "TJK", 2, 3, 1, -1, -1),
    @(2000, "TJK", 2, 3, 2, -1, -1),
    @(2001, "TJK", 2, 3, 3, -1, -1),
    @(2002, "TJK", 2, 3, 4, -1, -1),
    @(2003, "TJK", 1, 4, 5, -3, -3),
    @(2004, "TJK", 1, 4, 6, -3, -3),
    @(2005, "TJK", 1, 4, 7, -3, -3),
    @(2006, "TJK", 1, 4, 8, -3, -3),
    @(2007, "TJK", 1, 4, 9, -3, -3),
    @(2008, "TJK", 1, 4, 10, -3, -3),
    @(2009, "TJK", 1, 4, 11, -3, -3),
    @(2010, "TJK", 1, 4, 12, -3, -3),
    @(2011, "TJK", 1, 4, 13, -3, -3),
    @(2012, "TJK", 1, 4, 14, -3, -3),
    @(2013, "TJK", 1, 4, 15, -3, -3),
    @(2014, "TJK", 1, 4, 16, -3, -3),
    @(2015, "TJK", 1, 4, 17, -3, -3),
    @(2016, "TJK", 1, 4, 18, -3, -3),
    @(2017, "TJK", 1, 4, 19, -3, -3),
    @(2018, "TJK", 1, 4, 20, -3, -3),
    @(1991, "TKM", 0, 8, 0, -8, -8),
    @(1992, "TKM", 0, 9, 1, -9, -9),
    @(1993, "TKM", 0, 9, 2, -9, -9),
    @(1994, "TKM", 0, 9, 3, -9, -9),
    @(1995, "TKM", 0, 9, 4, -9, -9),
    @(1996, "TKM", 0, 9, 5, -9, -9),
    @(1997, "TKM", 0, 9, 6, -9, -9),
    @(1998, "TKM", 0, 9, 7, -9, -9),
    @(1999, "TKM", 0, 9, 8, -9, -9),
    @(2000, "TKM", 0, 9, 9, -9, -9),
    @(2001, "TKM", 0, 9, 10, -9, -9),
    @(2002, "TKM", 0, 9, 11, -9, -9),
    @(2003, "TKM", 0, 9, 12, -9, -9),
    @(2004, "TKM", 0, 9, 13, -9, -9),
    @(2005, "TKM", 0, 9, 14, -9, -9),
    @(2006, "TKM", 0, 9, 15, -9, -9),
    @(2007, "TKM", 0, 9, 16, -9, -9),
    @(2008, "TKM", 0, 9, 17, -9, -9),
    @(2009, "TKM", 0, 9, 18, -9, -9),
    @(2010, "TKM", 0, 9, 19, -9, -9),
    @(2011, "TKM", 0, 9, 20, -9, -9),
    @(2012, "TKM", 0, 9, 21, -9, -9),
    @(2013, "TKM", 0, 8, 22, -8, -8),
    @(2014, "TKM", 0, 8, 23, -8, -8),
    @(2015, "TKM", 0, 8, 24, -8, -8),
    @(2016, "TKM", 0, 8, 25, -8, -8),
    @(2017, "TKM", 0, 8, 26, -8, -8),
    @(2018, "TKM", 0, 8, 27, -8, -8),
    @(1991, "UKR", 6, 0, 0, 6, 6),
    @(1992, "UKR", 6, 0, 1, 6, 6),
    @(1993, "UKR", 5, 0, 2, 5, 5),
    @(1994, "UKR", 7, 0, 3, 7, 7),
    @(1995, "UKR", 7, 0, 4, 7, 7),
    @(1996, "UKR", 7, 0, 5, 7, 7),
    @(1997, "UKR", 7, 0, 6, 7, 7),
    @(1998, "UKR", 7, 0, 7, 7, 7),
    @(1999, "UKR", 7, 0, 8, 7, 7),
    @(2000, "UKR", 6, 0, 9, 6, 6),
    @(2001, "UKR", 6, 0, 10, 6, 6),
    @(2002, "UKR", 6, 0, 11, 6, 6),
    @(2003, "UKR", 6, 0, 12, 6, 6),
    @(2004, "UKR", 6, 0, 13, 6, 6),
    @(2005, "UKR", 6, 0, 14, 6, 6),
    @(2006, "UKR", 7, 0, 15, 7, 7),
    @(2007, "UKR", 7, 0, 16, 7, 7),
    @(2008, "UKR", 7, 0, 17, 7, 7),
    @(2009, "UKR", 7, 0, 18, 7, 7),
    @(2010, "UKR", 6, 0, 19, 6, 6),
    @(2011, "UKR", 6, 0, 20, 6, 6),
    @(2012, "UKR", 6, 0, 21, 6, 6),
    @(2013, "UKR", 6, 0, 22, 6, 6),
    @(2014, "UKR", 5, 1, 23, 4, 4),
    @(2015, "UKR", 5, 1, 24, 4, 4),
    @(2016, "UKR", 5, 1, 25, 4, 4),
    @(2017, "UKR", 5, 1, 26, 4, 4),
    @(2018, "UKR", 5, 1, 27, 4, 4),
    @(1991, "UZB", 0, 9, 0, -9, -9),
    @(1992, "UZB", 0, 9, 1, -9, -9),
    @(1993, "UZB", 0, 9, 2, -9, -9),
    @(1994, "UZB", 0, 9, 3, -9, -9),
    @(1995, "UZB", 0, 9, 4, -9, -9),
    @(1996, "UZB", 0, 9, 5, -9, -9),
    @(1997, "UZB", 0, 9, 6, -9, -9),
    @(1998, "UZB", 0, 9, 7, -9, -9),
    @(1999, "UZB", 0, 9, 8, -9, -9),
    @(2000, "UZB", 0, 9, 9, -9, -9),
    @(2001, "UZB", 0, 9, 10, -9, -9),
    @(2002, "UZB", 0, 9, 11, -9, -9),
    @(2003, "UZB", 0, 9, 12, -9, -9),
    @(2004, "UZB", 0, 9, 13, -9, -9),
    @(2005, "UZB", 0, 9, 14, -9, -9),
    @(2006, "UZB", 0, 9, 15, -9, -9),
    @(2007, "UZB", 0, 9, 16, -9, -9),
    @(2008, "UZB", 0, 9, 17, -9, -9),
    @(2009, "UZB", 0, 9, 18, -9, -9),
    @(2010, "UZB", 0, 9, 19, -9, -9),
    @(2011, "UZB", 0, 9, 20, -9, -9),
    @(2012, "UZB", 0, 9, 21, -9, -9),
    @(2013, "UZB", 0, 9, 22, -9, -9),
    @(2014, "UZB", 0, 9, 23, -9, -9),
    @(2015, "UZB", 0, 9, 24, -9, -9),
    @(2016, "UZB", 0, 9, 25, -9, -9),
    @(2017, "UZB", 0, 9, 26, -9, -9),
    @(2018, "UZB", 0, 9, 27, -9, -9)
)

$nRows = 125
$nCols = 7
$arr = New-Object 'object[,]' $nRows,$nCols
for ($i = 0; $i -lt $nRows; $i++) {
    for ($j = 0; $j -lt $nCols; $j++) {
        $arr[$i, $j] = $data[$i][$j]
    }
}

$rng = $ws.Range("A309:G433")
$rng.Value = $arr
